$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: extend the "${date_contract}" run with the extra clause naming
# the company and the employee placeholder. Find/Replace on the exact run
# text keeps the original run's rPr untouched, which matches the diff (the
# <w:rPr> block around the text is unchanged).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "`${date_contract}", $false, $false, $false, $false, $false, $true, 1, $false,
    "`${date_contract} giữa CÔNG TY TNHH DỊCH VỤ TRẤN THANH và Ông/Bà `${employee}", 2
) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: the paragraph right after it (previously an almost-empty
# paragraph) is rebuilt with new paragraph properties (tab stops, spacing,
# indent) and a single run containing a tab followed by "${base}".
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("`${date_contract} giữa CÔNG TY TNHH DỊCH VỤ TRẤN THANH và Ông/Bà `${employee}") | Out-Null
$datePara = $rng.Paragraphs(1)
$baseParaRange = $datePara.Next().Range

$baseXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="LOnormal"/><w:tabs><w:tab w:val="clear" w:pos="720"/><w:tab w:val="left" w:pos="360" w:leader="none"/></w:tabs><w:spacing w:lineRule="auto" w:line="360" w:before="120" w:after="120"/><w:ind w:hanging="0"/><w:jc w:val="both"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:position w:val="0"/><w:sz w:val="20"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:highlight w:val="white"/><w:vertAlign w:val="baseline"/></w:rPr><w:tab/><w:t>${base}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$baseParaRange.InsertXML($baseXml) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: the address paragraph splits its single run into two runs -
# "Địa chỉ: " (preserving the trailing space) and the new address text
# ("khu phố Scenic Valley 2, " dropped, "Nguyễn Văn Linh" becomes
# "Đường Nguyễn Văn Linh").
# ---------------------------------------------------------------------------
$addrRng = $d.Content
$addrRng.Find.Execute("Địa chỉ: A0.01 khu phố Scenic Valley 2, Nguyễn Văn Linh, Phường Tân Phú, Quận 7, TPHCM") | Out-Null
$addrPara = $addrRng.Paragraphs(1)

$addrXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="LOnormal"/><w:spacing w:lineRule="auto" w:line="360"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Địa chỉ: </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>A0.01 Đường Nguyễn Văn Linh, Phường Tân Phú, Quận 7, TPHCM</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$addrPara.Range.InsertXML($addrXml) | Out-Null

Write-Host "edits applied"
